$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# Update End Year from 2050 to 2024
$ws.Range("B4").Value = 2024

# Add formula to start_year_dismantling (B12) = B7 + 1 (value remains 4)
$ws.Range("B12").Formula = "=B7+1"

# Reduce maximum_investment_capacity_per_year from 1000000 to 100
$ws.Range("B13").Value = 100

# Flip realistic_candidate_capacities_tobe_installed (B16) to TRUE
$ws.Range("B16").Value = $true

# Flip realistic_candidate_capacities_for_future (B17) to TRUE
$ws.Range("B17").Value = $true

# Reduce dummy_capacity (B19) from 300 to 100
$ws.Range("B19").Value = 100

# Update selection to C8
$ws.Range("C8").Select()
